$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and row 3 hold duplicate player/team info but their runs/balls/fours
# stats (columns C:E) were swapped between the two innings rows. Use
# Copy/Paste (rather than direct .Value assignment) so the cells keep their
# original "stored as text" type instead of being coerced to numbers.

$buffer = $ws.Range("C100:E100")

$ws.Range("C2:E2").Copy($buffer)
$ws.Range("C3:E3").Copy($ws.Range("C2:E2"))
$buffer.Copy($ws.Range("C3:E3"))

$buffer.Clear()
